$wb = $excel.ActiveWorkbook
$todos = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add($null, $todos)
$newSheet.Name = "Users"

$newSheet.Range("A1").Value = "Id"
$newSheet.Range("B1").Value = "Username"
$newSheet.Range("D1").Value = "Firstname"
$newSheet.Range("E1").Value = "Lastname"
$newSheet.Range("F1").Value = "Email"
$newSheet.Range("G1").Value = "IsDeleted"
$newSheet.Range("C1").Value = "Password"

[void]$newSheet.Range("C2").Select()
